# Actualizacion automatica: inserta un nuevo cliente "PROVEEDORA PARA
# METALMECANICA E INDUSTRIAS PROMETIN CIA LTDA" antes de "QUINTEROS VELASCO
# ELSA ROSARIO" en las hojas "VENTAS POR GRUPO" y "VENTA MENSUAL", y ajusta
# el ancho de la columna B y los rotulos de totales en consecuencia.

$wb = $excel.ActiveWorkbook

$newClientName = "PROVEEDORA PARA METALMECANICA E INDUSTRIAS PROMETIN CIA LTDA"
$asesorName = "ALMEIDA CUATIN JHONATHANN CARLOS"

# Ancho exacto deseado para la columna B (62 caracteres). Excel agrega un
# relleno interno de 5/6 de caracter al guardar el ancho, asi que restamos
# esa cantidad al asignar ColumnWidth para que el ancho resultante en el XML
# sea exactamente 62.
$targetColBWidth = 62
$colWidthPadding = 5 / 6
$colBWidthToSet = $targetColBWidth - $colWidthPadding

# ---------------------------------------------------------------------
# Hoja "VENTAS POR GRUPO": columnas A..R, fila de insercion = 26
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("VENTAS POR GRUPO")

$ws1.Rows(26).Insert()

$ws1.Range("A26").Value = $asesorName
$ws1.Range("B26").Value = $newClientName

$cols1 = @("C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R")
foreach ($col in $cols1) {
    $ws1.Range($col + "26").Value = 0
}

$ws1.Columns("B").ColumnWidth = $colBWidthToSet

# Actualiza los rotulos "X de 33" -> "X de 34" en la nueva fila de totales (36)
foreach ($col in $cols1) {
    $cell = $ws1.Range($col + "36")
    $cell.Value = $cell.Text -replace "de 33", "de 34"
}

# ---------------------------------------------------------------------
# Hoja "VENTA MENSUAL": columnas A..G, fila de insercion = 26
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("VENTA MENSUAL")

$ws2.Rows(26).Insert()

$ws2.Range("A26").Value = $asesorName
$ws2.Range("B26").Value = $newClientName

$cols2 = @("C","D","E","F","G")
foreach ($col in $cols2) {
    $ws2.Range($col + "26").Value = 0
}

$ws2.Columns("B").ColumnWidth = $colBWidthToSet
